# Restore cell C10 ("From" value for rule R30) on the Rules sheet
# from 18 to 1, as captured by the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
